$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2402.8
$ws.Range("I15").Value = 2402.8
$ws.Range("K15").Value = 7208.400000000001
$ws.Range("M15").Value = -7039.400000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6981.905
$ws.Range("I116").Value = 8552.5
$ws.Range("K116").Value = 8552.5
$ws.Range("M116").Value = -5110.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2639.8
$ws.Range("I132").Value = 2480.762
$ws.Range("K132").Value = 7442.286
$ws.Range("M132").Value = -4912.286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2249.5715
$ws.Range("I137").Value = 2258
$ws.Range("K137").Value = 6774
$ws.Range("M137").Value = -4224

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2091.8442
$ws.Range("I138").Value = 2276.8572
$ws.Range("J138").Value = 2022.4642
$ws.Range("K138").Value = 6830.571599999999
$ws.Range("L138").Value = 6067.392599999999
$ws.Range("M138").Value = -1690.571599999999
$ws.Range("N138").Value = -16347.3926

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1788.25
$ws.Range("I74").Value = 1232.5385
$ws.Range("K74").Value = 1232.5385
$ws.Range("M74").Value = -358.5385000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1788.25
$ws.Range("I77").Value = 1232.5385
$ws.Range("K77").Value = 6162.692500000001
$ws.Range("M77").Value = -1794.692500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 25428
$ws.Range("J123").Value = 25428
$ws.Range("L123").Value = 25428
$ws.Range("N123").Value = -35228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 400.64706
$ws.Range("I64").Value = 431.54544
$ws.Range("J64").Value = 344
$ws.Range("K64").Value = 431.54544
$ws.Range("L64").Value = 344
$ws.Range("M64").Value = -206.54544
$ws.Range("N64").Value = -794

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 400.64706
$ws.Range("I67").Value = 431.54544
$ws.Range("J67").Value = 344
$ws.Range("K67").Value = 431.54544
$ws.Range("L67").Value = 344
$ws.Range("M67").Value = 348.45456
$ws.Range("N67").Value = -1904

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2498.9167
$ws.Range("J134").Value = 3193.7646
$ws.Range("L134").Value = 9581.293799999999
$ws.Range("N134").Value = -14651.2938

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5626.9375
$ws.Range("I31").Value = 1223.0416
$ws.Range("J31").Value = 10030.833
$ws.Range("K31").Value = 1223.0416
$ws.Range("L31").Value = 10030.833
$ws.Range("M31").Value = -928.0416
$ws.Range("N31").Value = -10620.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5626.9375
$ws.Range("I34").Value = 1223.0416
$ws.Range("J34").Value = 10030.833
$ws.Range("K34").Value = 1223.0416
$ws.Range("L34").Value = 10030.833
$ws.Range("M34").Value = -1021.0416
$ws.Range("N34").Value = -10434.833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1175.5238
$ws.Range("I58").Value = 940.2174
$ws.Range("J58").Value = 1460.3684
$ws.Range("K58").Value = 940.2174
$ws.Range("L58").Value = 1460.3684
$ws.Range("M58").Value = -737.2174
$ws.Range("N58").Value = -1866.3684

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1986.3636
$ws.Range("I99").Value = 1850
$ws.Range("K99").Value = 1850
$ws.Range("M99").Value = -352

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2022.2
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1986.3636
$ws.Range("I126").Value = 1850
$ws.Range("K126").Value = 5550
$ws.Range("M126").Value = -3080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3300.2727
$ws.Range("I132").Value = 3361.2
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 10083.6
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -7553.599999999999
$ws.Range("N132").Value = -14808.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9656.615
$ws.Range("I134").Value = 11653.6
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 34960.8
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -32425.8
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1175.5238
$ws.Range("I136").Value = 940.2174
$ws.Range("J136").Value = 1460.3684
$ws.Range("K136").Value = 2820.6522
$ws.Range("L136").Value = 4381.1052
$ws.Range("M136").Value = -270.6522
$ws.Range("N136").Value = -9481.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 75000136
$ws.Range("I23").Value = 142.5
$ws.Range("J23").Value = 93750136
$ws.Range("K23").Value = 427.5
$ws.Range("L23").Value = 281250408
$ws.Range("M23").Value = -192.5
$ws.Range("N23").Value = -281250878

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 5618.2856
$ws.Range("J123").Value = 9999.333000000001
$ws.Range("L123").Value = 29997.999
$ws.Range("N123").Value = -34897.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2476.8157
$ws.Range("I132").Value = 2054.923
$ws.Range("J132").Value = 2696.2
$ws.Range("K132").Value = 18494.307
$ws.Range("L132").Value = 24265.8
$ws.Range("M132").Value = -15964.307
$ws.Range("N132").Value = -29325.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2321.8333
$ws.Range("I31").Value = 2026.2
$ws.Range("J31").Value = 3800
$ws.Range("K31").Value = 2026.2
$ws.Range("L31").Value = 3800
$ws.Range("M31").Value = -1734.2
$ws.Range("N31").Value = -4384

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 2321.8333
$ws.Range("I37").Value = 2026.2
$ws.Range("J37").Value = 3800
$ws.Range("K37").Value = 2026.2
$ws.Range("L37").Value = 3800
$ws.Range("M37").Value = -1749.2
$ws.Range("N37").Value = -4354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 27428.715
$ws.Range("I22").Value = 2399.6
$ws.Range("K22").Value = 2399.6
$ws.Range("M22").Value = -2104.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 27428.715
$ws.Range("I27").Value = 2399.6
$ws.Range("K27").Value = 2399.6
$ws.Range("M27").Value = -2292.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3221.4473
$ws.Range("I132").Value = 3033.625
$ws.Range("K132").Value = 9100.875
$ws.Range("M132").Value = -6570.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5052251.5
$ws.Range("I136").Value = 1700.6428
$ws.Range("J136").Value = 33335336
$ws.Range("K136").Value = 5101.928400000001
$ws.Range("L136").Value = 100006008
$ws.Range("M136").Value = -2551.928400000001
$ws.Range("N136").Value = -100011108

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 22357.5
$ws.Range("I5").Value = 4999
$ws.Range("K5").Value = 4999
$ws.Range("M5").Value = -4887

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2112.7568
$ws.Range("I122").Value = 2030.6875
$ws.Range("J122").Value = 2638
$ws.Range("K122").Value = 6092.0625
$ws.Range("L122").Value = 7914
$ws.Range("M122").Value = -3642.0625
$ws.Range("N122").Value = -12814

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 32047.666
$ws.Range("J123").Value = 56143
$ws.Range("L123").Value = 56143
$ws.Range("N123").Value = -65943

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1928.5063
$ws.Range("I136").Value = 1627.4762
$ws.Range("J136").Value = 3113.8125
$ws.Range("K136").Value = 4882.4286
$ws.Range("L136").Value = 9341.4375
$ws.Range("M136").Value = -2332.4286
$ws.Range("N136").Value = -14441.4375
